$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-21 for columns I (I0) and J (IF)
$data = @(
    @(1, 4),
    @(1, 6),
    @(1, 3),
    @(1, 6),
    @(1, 3),
    @(1, 5),
    @(1, 6),
    @(1, 7),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 8),
    @(1, 3),
    @(1, 1),
    @(1, 5),
    @(1, 2),
    @(5, 7),
    @(1, 3),
    @(5, 6),
    @(3, 4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
